$p = $ppt.ActivePresentation

# --- 1. Update the fixed "datetimeFigureOut" date text from 17-08-2020 to 18-08-2020
#        everywhere it appears: the slide master, every slide layout, and the notes master.

$oldDate = "17-08-2020"
$newDate = "18-08-2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Add the two new text boxes (author name + student id) to slide 1.
#        Positions/sizes below are in points (EMU / 12700), matching the
#        target offsets 8644919,213064 / 1884106x369332 EMU and
#        9586972,610794 / 1237839x369332 EMU.

$s1 = $p.Slides.Item(1)

$tb1 = $s1.Shapes.AddTextbox(1, 680.7022834645669, 16.776692913385826, 148.3548031496063, 29.081259842519685)
$tb1.Name = "Tekstfelt 3"
$tb1.Fill.Visible = $false
$tb1.TextFrame.WordWrap = $false
$tb1.TextFrame.AutoSize = 1
$tb1.TextFrame.TextRange.Text = "Andreas Blaabjerg"
$tb1.TextFrame.TextRange.LanguageID = "da-DK"

$tb2 = $s1.Shapes.AddTextbox(1, 754.8796850393701, 48.094015748031495, 97.46763779527559, 29.081259842519685)
$tb2.Name = "Tekstfelt 4"
$tb2.Fill.Visible = $false
$tb2.TextFrame.WordWrap = $false
$tb2.TextFrame.AutoSize = 1
$tb2.TextFrame.TextRange.Text = "201510924"
$tb2.TextFrame.TextRange.LanguageID = "da-DK"
